# Apply the crypto-tracker refresh update to Sheet1.
# For each changed row we update Price (D) / Volume(1h) (E), and for the
# four rows whose coins were reordered we also update Coin (B) and Link (C).
# Numeric-looking Price strings are written with a leading apostrophe so Excel
# keeps storing them as text (matching the original "dotted" text formatting)
# instead of silently converting them into floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '68.675.56'
$ws.Cells.Item(2, 5).Value = '  +0.66%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.707.74'
$ws.Cells.Item(3, 5).Value = '  +2.27%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.04%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''599.41'
$ws.Cells.Item(5, 5).Value = '  +0.31%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''160.93'
$ws.Cells.Item(6, 5).Value = '  +2.67%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.02%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.544'
$ws.Cells.Item(8, 5).Value = '  +0.11%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.706.51'
$ws.Cells.Item(9, 5).Value = '  +2.26%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +0.29%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.32%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +1.00%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +2.86%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''28.40'
$ws.Cells.Item(14, 5).Value = '  +1.28%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.194.00'
$ws.Cells.Item(15, 5).Value = '  +2.07%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''0.0000189'
$ws.Cells.Item(16, 5).Value = '  -0.78%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '68.569.28'
$ws.Cells.Item(17, 5).Value = '  +0.56%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.705.33'
$ws.Cells.Item(18, 5).Value = '  +1.76%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +4.32%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''7.68'
$ws.Cells.Item(20, 5).Value = '  +4.59%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''365.04'
$ws.Cells.Item(21, 5).Value = '  +0.47%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +2.87%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''4.91'
$ws.Cells.Item(23, 5).Value = '  +2.18%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +2.33%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''74.17'
$ws.Cells.Item(25, 5).Value = '  -1.40%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.06%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''9.96'
$ws.Cells.Item(27, 5).Value = '  +2.06%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '2.838.60'
$ws.Cells.Item(28, 5).Value = '  +2.12%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''0.0000106'
$ws.Cells.Item(29, 5).Value = '  +1.01%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''594.81'
$ws.Cells.Item(30, 5).Value = '  +6.44%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -0.01%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''8.27'
$ws.Cells.Item(32, 5).Value = '  +2.59%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +2.86%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +4.70%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +2.97%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +5.46%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.06%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''161.44'
$ws.Cells.Item(38, 5).Value = '  +0.28%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''19.86'

# Row 40
$ws.Cells.Item(40, 5).Value = '  +2.31%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +2.21%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''5.40'
$ws.Cells.Item(42, 5).Value = '  +1.48%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'dogwifhat'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(43, 4).Value = '''2.68'
$ws.Cells.Item(43, 5).Value = '  +2.92%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(44, 4).Value = '''18.01'
$ws.Cells.Item(44, 5).Value = '  +1.20%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(45, 4).Value = '0.0₆0317'
$ws.Cells.Item(45, 5).Value = '  -5.65%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'USDe'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(46, 4).Value = '''1.00'
$ws.Cells.Item(46, 5).Value = '  +0.03%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''158.14'
$ws.Cells.Item(47, 5).Value = '  -0.59%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +5.71%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''1.78'
$ws.Cells.Item(49, 5).Value = '  +5.29%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''0.606'

# Row 51
$ws.Cells.Item(51, 4).Value = '''22.13'
$ws.Cells.Item(51, 5).Value = '  +0.25%  '
